$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new "Save" column, matching the style used by the other headers (G1 etc.)
$ws.Range("H1").Value = "Save"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = $ws.Range("G1").HorizontalAlignment
$ws.Range("H1").VerticalAlignment = $ws.Range("G1").VerticalAlignment
$ws.Range("H1").Borders.LineStyle = $ws.Range("G1").Borders.LineStyle

$saveValues = @(0,1,0,0,0,0,0,0,0,0,0,0,0,1,0,1,0,0,0,1,0,1,1,1,0,0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
